# Intraday standard output (Mar-Sab)
# Insert the two new ticker rows, then populate rows 2-7 with the refreshed dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AHMA becomes the new row 2 (ATON and the rows below it shift down by one)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# BCTX becomes the new row 4 (DCOY, EVTV, XAIR shift down by one more)
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).ClearFormats()

# Row 2: AHMA
$ws.Range("A2").Value = 'AHMA'
$ws.Range("B2").NumberFormat = "YYYY-MM-DD"
$ws.Range("B2").Value = "2026-01-13"
$ws.Range("C2").Value = 113.83
$ws.Range("D2").Value = 364030000
$ws.Range("E2").Value = 1920000
$ws.Range("F2").Value = 10740000
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '9.30%'
$ws.Range("G2").NumberFormat = "General"
$ws.Range("H2").Value = 93.49
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = 0.41
$ws.Range("K2").Value = 10.96
$ws.Range("L2").Value = 11.22
$ws.Range("M2").Value = 13.83
$ws.Range("N2").Value = 8.82
$ws.Range("O2").Value = 12.31
$ws.Range("P2").Value = 83662202
$ws.Range("Q2").Value = '2026-01-13 14:59:00'
$ws.Range("R2").Value = '2026-01-13 09:58:00'
$ws.Range("S2").Value = 5.15
$ws.Range("T2").Value = 15.21
$ws.Range("U2").Value = 5.07
$ws.Range("V2").Value = 11.3
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = '2026-01-13 08:22:00'
$ws.Range("Y2").Value = 11.29
$ws.Range("Z2").Value = 10.68
$ws.Range("AA2").Value = 22797890
$ws.Range("AB2").Value = 12.12
$ws.Range("AC2").Value = 10.25
$ws.Range("AD2").Value = 23309466
$ws.Range("AE2").Value = 12.5
$ws.Range("AF2").Value = 8.82
$ws.Range("AG2").Value = 38114703
$ws.Range("AH2").Value = 12.5
$ws.Range("AI2").Value = 8.82
$ws.Range("AJ2").Value = 40297229
$ws.Range("AK2").Value = 10.38
$ws.Range("AL2").Value = 12.5
$ws.Range("AM2").Value = 8.82
$ws.Range("AN2").Value = 41373362
$ws.Range("AO2").Value = 10.21
$ws.Range("AP2").Value = 12.5
$ws.Range("AQ2").Value = 8.82
$ws.Range("AR2").Value = 42011162
$ws.Range("AS2").Value = 11.23
$ws.Range("AT2").Value = 12.5
$ws.Range("AU2").Value = 8.82
$ws.Range("AV2").Value = 61492984
$ws.Range("AW2").Value = 12.33

# Row 3: ATON
$ws.Range("A3").Value = 'ATON'
$ws.Range("B3").NumberFormat = "YYYY-MM-DD"
$ws.Range("B3").Value = "2026-01-13"
$ws.Range("C3").Value = 102.2
$ws.Range("D3").Value = 9910000
$ws.Range("E3").Value = 2820000
$ws.Range("F3").Value = 5740000
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '-29.35%'
$ws.Range("G3").NumberFormat = "General"
$ws.Range("H3").Value = 62.99
$ws.Range("I3").Value = 0.52
$ws.Range("J3").Value = 8.78
$ws.Range("K3").Value = 1.88
$ws.Range("L3").Value = 1.82
$ws.Range("M3").Value = 3.3
$ws.Range("N3").Value = 1.53
$ws.Range("O3").Value = 1.99
$ws.Range("P3").Value = 385413301
$ws.Range("Q3").Value = '2026-01-13 15:18:00'
$ws.Range("R3").Value = '2026-01-13 10:06:00'
$ws.Range("S3").Value = 2.59
$ws.Range("T3").Value = 2.96
$ws.Range("U3").Value = 1.81
$ws.Range("V3").Value = 1.81
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = '2026-01-13 08:04:00'
$ws.Range("Y3").Value = 1.93
$ws.Range("Z3").Value = 1.77
$ws.Range("AA3").Value = 30969502
$ws.Range("AB3").Value = 2.01
$ws.Range("AC3").Value = 1.77
$ws.Range("AD3").Value = 32480567
$ws.Range("AE3").Value = 2.19
$ws.Range("AF3").Value = 1.72
$ws.Range("AG3").Value = 42146099
$ws.Range("AH3").Value = 2.19
$ws.Range("AI3").Value = 1.53
$ws.Range("AJ3").Value = 49518119
$ws.Range("AK3").Value = 1.81
$ws.Range("AL3").Value = 2.19
$ws.Range("AM3").Value = 1.53
$ws.Range("AN3").Value = 130307429
$ws.Range("AO3").Value = 1.73
$ws.Range("AP3").Value = 2.19
$ws.Range("AQ3").Value = 1.53
$ws.Range("AR3").Value = 175070886
$ws.Range("AS3").Value = 1.91
$ws.Range("AT3").Value = 2.19
$ws.Range("AU3").Value = 1.53
$ws.Range("AV3").Value = 231107242
$ws.Range("AW3").Value = 2.01

# Row 4: BCTX
$ws.Range("A4").Value = 'BCTX'
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"
$ws.Range("B4").Value = "2026-01-13"
$ws.Range("C4").Value = 57.07
$ws.Range("D4").Value = 20530000
$ws.Range("E4").Value = 1860000
$ws.Range("F4").Value = 1880000
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '-9.00%'
$ws.Range("G4").NumberFormat = "General"
$ws.Range("H4").Value = 1.27
$ws.Range("I4").Value = 1.67
$ws.Range("J4").Value = 15.13
$ws.Range("K4").Value = 11.65
$ws.Range("L4").Value = 11.96
$ws.Range("M4").Value = 12.1
$ws.Range("N4").Value = 9.06
$ws.Range("O4").Value = 11.15
$ws.Range("P4").Value = 12582674
$ws.Range("Q4").Value = '2026-01-13 09:30:00'
$ws.Range("R4").Value = '2026-01-13 09:56:00'
$ws.Range("S4").Value = 8
$ws.Range("T4").Value = 19.68
$ws.Range("U4").Value = 7.82
$ws.Range("V4").Value = 12.27
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = '2026-01-13 08:12:00'
$ws.Range("Y4").Value = 12.1
$ws.Range("Z4").Value = 11.3
$ws.Range("AA4").Value = 3972208
$ws.Range("AB4").Value = 12.1
$ws.Range("AC4").Value = 10.41
$ws.Range("AD4").Value = 4027085
$ws.Range("AE4").Value = 12.1
$ws.Range("AF4").Value = 9.06
$ws.Range("AG4").Value = 4523387
$ws.Range("AH4").Value = 12.1
$ws.Range("AI4").Value = 9.06
$ws.Range("AJ4").Value = 4747625
$ws.Range("AK4").Value = 10.2
$ws.Range("AL4").Value = 12.1
$ws.Range("AM4").Value = 9.06
$ws.Range("AN4").Value = 4948881
$ws.Range("AO4").Value = 10.52
$ws.Range("AP4").Value = 12.1
$ws.Range("AQ4").Value = 9.06
$ws.Range("AR4").Value = 5081171
$ws.Range("AS4").Value = 10.12
$ws.Range("AT4").Value = 12.1
$ws.Range("AU4").Value = 9.06
$ws.Range("AV4").Value = 5271397
$ws.Range("AW4").Value = 10.4

# Row 5: DCOY
$ws.Range("A5").Value = 'DCOY'
$ws.Range("B5").NumberFormat = "YYYY-MM-DD"
$ws.Range("B5").Value = "2026-01-13"
$ws.Range("C5").Value = 65.64
$ws.Range("D5").Value = 7980000
$ws.Range("E5").Value = 6000000
$ws.Range("F5").Value = 6380000
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '-7.41%'
$ws.Range("G5").NumberFormat = "General"
$ws.Range("H5").Value = 6.05
$ws.Range("I5").Value = 4.71
$ws.Range("J5").Value = 3.08
$ws.Range("K5").Value = 1.32
$ws.Range("L5").Value = 1.35
$ws.Range("M5").Value = 1.41
$ws.Range("N5").Value = 1.07
$ws.Range("O5").Value = 1.23
$ws.Range("P5").Value = 825232430
$ws.Range("Q5").Value = '2026-01-13 09:30:00'
$ws.Range("R5").Value = '2026-01-13 09:53:00'
$ws.Range("S5").Value = 0.8
$ws.Range("T5").Value = 1.81
$ws.Range("U5").Value = 0.77
$ws.Range("V5").Value = 1.36
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = '2026-01-13 08:43:00'
$ws.Range("Y5").Value = 1.41
$ws.Range("Z5").Value = 1.26
$ws.Range("AA5").Value = 51828392
$ws.Range("AB5").Value = 1.41
$ws.Range("AC5").Value = 1.17
$ws.Range("AD5").Value = 54446957
$ws.Range("AE5").Value = 1.41
$ws.Range("AF5").Value = 1.07
$ws.Range("AG5").Value = 62081185
$ws.Range("AH5").Value = 1.41
$ws.Range("AI5").Value = 1.07
$ws.Range("AJ5").Value = 184284311
$ws.Range("AK5").Value = 1.23
$ws.Range("AL5").Value = 1.41
$ws.Range("AM5").Value = 1.07
$ws.Range("AN5").Value = 227785782
$ws.Range("AO5").Value = 1.14
$ws.Range("AP5").Value = 1.41
$ws.Range("AQ5").Value = 1.07
$ws.Range("AR5").Value = 229893959
$ws.Range("AS5").Value = 1.17
$ws.Range("AT5").Value = 1.41
$ws.Range("AU5").Value = 1.07
$ws.Range("AV5").Value = 361635477
$ws.Range("AW5").Value = 1.22

# Row 6: EVTV
$ws.Range("A6").Value = 'EVTV'
$ws.Range("B6").NumberFormat = "YYYY-MM-DD"
$ws.Range("B6").Value = "2026-01-13"
$ws.Range("C6").Value = 33.07
$ws.Range("D6").Value = 17000000
$ws.Range("E6").Value = 4450000
$ws.Range("F6").Value = 4830000
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '5.39%'
$ws.Range("G6").NumberFormat = "General"
$ws.Range("H6").Value = 7.9
$ws.Range("I6").Value = 1.91
$ws.Range("J6").Value = 0.97
$ws.Range("K6").Value = 3.21
$ws.Range("L6").Value = 3.32
$ws.Range("M6").Value = 4.91
$ws.Range("N6").Value = 2.62
$ws.Range("O6").Value = 3.52
$ws.Range("P6").Value = 422796722
$ws.Range("Q6").Value = '2026-01-13 13:34:00'
$ws.Range("R6").Value = '2026-01-13 09:54:00'
$ws.Range("S6").Value = 3.73
$ws.Range("T6").Value = 3.95
$ws.Range("U6").Value = 2.88
$ws.Range("V6").Value = 3.34
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = '2026-01-13 08:02:00'
$ws.Range("Y6").Value = 3.34
$ws.Range("Z6").Value = 3.08
$ws.Range("AA6").Value = 30901272
$ws.Range("AB6").Value = 3.57
$ws.Range("AC6").Value = 2.95
$ws.Range("AD6").Value = 37796996
$ws.Range("AE6").Value = 3.57
$ws.Range("AF6").Value = 2.62
$ws.Range("AG6").Value = 51194768
$ws.Range("AH6").Value = 4.02
$ws.Range("AI6").Value = 2.62
$ws.Range("AJ6").Value = 77204596
$ws.Range("AK6").Value = 3.57
$ws.Range("AL6").Value = 4.02
$ws.Range("AM6").Value = 2.62
$ws.Range("AN6").Value = 91473708
$ws.Range("AO6").Value = 3.67
$ws.Range("AP6").Value = 4.02
$ws.Range("AQ6").Value = 2.62
$ws.Range("AR6").Value = 162732531
$ws.Range("AS6").Value = 3.45
$ws.Range("AT6").Value = 4.87
$ws.Range("AU6").Value = 2.62
$ws.Range("AV6").Value = 284893748
$ws.Range("AW6").Value = 4.51

# Row 7: XAIR
$ws.Range("A7").Value = 'XAIR'
$ws.Range("B7").NumberFormat = "YYYY-MM-DD"
$ws.Range("B7").Value = "2026-01-13"
$ws.Range("C7").Value = 79.22
$ws.Range("D7").Value = 17550000
$ws.Range("E7").Value = 7540000
$ws.Range("F7").Value = 8010000
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '39.49%'
$ws.Range("G7").NumberFormat = "General"
$ws.Range("H7").Value = 5.87
$ws.Range("I7").Value = 18.88
$ws.Range("J7").Value = 3.14
$ws.Range("K7").Value = 1.63
$ws.Range("L7").Value = 1.57
$ws.Range("M7").Value = 2.66
$ws.Range("N7").Value = 1.37
$ws.Range("O7").Value = 2.17
$ws.Range("P7").Value = 1396493630
$ws.Range("Q7").Value = '2026-01-13 12:00:00'
$ws.Range("R7").Value = '2026-01-13 09:35:00'
$ws.Range("S7").Value = 0.89
$ws.Range("T7").Value = 1.75
$ws.Range("U7").Value = 0.89
$ws.Range("V7").Value = 1.64
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = '2026-01-13 09:19:00'
$ws.Range("Y7").Value = 1.69
$ws.Range("Z7").Value = 1.55
$ws.Range("AA7").Value = 60403890
$ws.Range("AB7").Value = 1.69
$ws.Range("AC7").Value = 1.55
$ws.Range("AD7").Value = 64815150
$ws.Range("AE7").Value = 2.04
$ws.Range("AF7").Value = 1.37
$ws.Range("AG7").Value = 283855847
$ws.Range("AH7").Value = 2.5
$ws.Range("AI7").Value = 1.37
$ws.Range("AJ7").Value = 399764564
$ws.Range("AK7").Value = 2.38
$ws.Range("AL7").Value = 2.5
$ws.Range("AM7").Value = 1.37
$ws.Range("AN7").Value = 528265657
$ws.Range("AO7").Value = 2.2
$ws.Range("AP7").Value = 2.5
$ws.Range("AQ7").Value = 1.37
$ws.Range("AR7").Value = 540115272
$ws.Range("AS7").Value = 2.51
$ws.Range("AT7").Value = 2.66
$ws.Range("AU7").Value = 1.37
$ws.Range("AV7").Value = 878796407
$ws.Range("AW7").Value = 2.34

